{"js": "// Replace the date line and the 25 two-digit multiplication problems with\n// their updated values. Every \"old\" string below is unique in the\n// document, so a direct search-and-replace for each pair is unambiguous.\nconst replacements = [\n  [\"2025-05-10 Saturday\", \"2025-05-11 Sunday\"],\n  [\"47\\u00D788=\", \"94\\u00D724=\"],\n  [\"72\\u00D740=\", \"32\\u00D753=\"],\n  [\"86\\u00D754=\", \"92\\u00D731=\"],\n  [\"99\\u00D789=\", \"25\\u00D750=\"],\n  [\"68\\u00D763=\", \"28\\u00D744=\"],\n  [\"66\\u00D759=\", \"64\\u00D725=\"],\n  [\"18\\u00D740=\", \"93\\u00D750=\"],\n  [\"77\\u00D732=\", \"19\\u00D728=\"],\n  [\"81\\u00D719=\", \"31\\u00D756=\"],\n  [\"59\\u00D757=\", \"68\\u00D737=\"],\n  [\"24\\u00D716=\", \"39\\u00D719=\"],\n  [\"49\\u00D718=\", \"23\\u00D714=\"],\n  [\"86\\u00D724=\", \"32\\u00D765=\"],\n  [\"87\\u00D711=\", \"86\\u00D771=\"],\n  [\"89\\u00D796=\", \"33\\u00D776=\"],\n  [\"21\\u00D791=\", \"15\\u00D727=\"],\n  [\"67\\u00D741=\", \"52\\u00D762=\"],\n  [\"45\\u00D788=\", \"47\\u00D728=\"],\n  [\"46\\u00D732=\", \"93\\u00D757=\"],\n  [\"35\\u00D774=\", \"85\\u00D783=\"],\n  [\"99\\u00D788=\", \"55\\u00D757=\"],\n  [\"39\\u00D762=\", \"35\\u00D736=\"],\n  [\"14\\u00D762=\", \"34\\u00D746=\"],\n  [\"83\\u00D758=\", \"44\\u00D769=\"],\n  [\"95\\u00D765=\", \"46\\u00D713=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit multiplication problems with\n# their updated values. Every \"old\" string below is unique in the\n# document, so a Find/Replace (wdReplaceAll) for each pair is unambiguous.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    @(\"2025-05-10 Saturday\", \"2025-05-11 Sunday\"),\n    @(\"47\u00d788=\", \"94\u00d724=\"),\n    @(\"72\u00d740=\", \"32\u00d753=\"),\n    @(\"86\u00d754=\", \"92\u00d731=\"),\n    @(\"99\u00d789=\", \"25\u00d750=\"),\n    @(\"68\u00d763=\", \"28\u00d744=\"),\n    @(\"66\u00d759=\", \"64\u00d725=\"),\n    @(\"18\u00d740=\", \"93\u00d750=\"),\n    @(\"77\u00d732=\", \"19\u00d728=\"),\n    @(\"81\u00d719=\", \"31\u00d756=\"),\n    @(\"59\u00d757=\", \"68\u00d737=\"),\n    @(\"24\u00d716=\", \"39\u00d719=\"),\n    @(\"49\u00d718=\", \"23\u00d714=\"),\n    @(\"86\u00d724=\", \"32\u00d765=\"),\n    @(\"87\u00d711=\", \"86\u00d771=\"),\n    @(\"89\u00d796=\", \"33\u00d776=\"),\n    @(\"21\u00d791=\", \"15\u00d727=\"),\n    @(\"67\u00d741=\", \"52\u00d762=\"),\n    @(\"45\u00d788=\", \"47\u00d728=\"),\n    @(\"46\u00d732=\", \"93\u00d757=\"),\n    @(\"35\u00d774=\", \"85\u00d783=\"),\n    @(\"99\u00d788=\", \"55\u00d757=\"),\n    @(\"39\u00d762=\", \"35\u00d736=\"),\n    @(\"14\u00d762=\", \"34\u00d746=\"),\n    @(\"83\u00d758=\", \"44\u00d769=\"),\n    @(\"95\u00d765=\", \"46\u00d713=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n}\n"}
